$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1363
$ws1.Range("F3").Value = 2893
$ws1.Range("F4").Value = 5

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1363
$ws4.Range("F4").Value = 2893
$ws4.Range("F5").Value = 5
